$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: fix G3 from -0 to 0
$ws.Range("G3").Value = 0

# Row 4
$ws.Range("E4").Value = 5
$ws.Range("G4").Value = -3
$ws.Range("H4").Value = 13

# Row 8
$ws.Range("E8").Value = 5
$ws.Range("G8").Value = -3
$ws.Range("H8").Value = 13

# Row 11: fix G11 from -0 to 0
$ws.Range("G11").Value = 0

# Row 15: fix G15 from -0 to 0
$ws.Range("G15").Value = 0

# Row 16
$ws.Range("E16").Value = 6
$ws.Range("G16").Value = -3
$ws.Range("H16").Value = 13

# Row 18
$ws.Range("E18").Value = 5
$ws.Range("G18").Value = -3
$ws.Range("H18").Value = 13

# Row 21: fix G21 from -0 to 0
$ws.Range("G21").Value = 0

# Row 23
$ws.Range("E23").Value = 4
$ws.Range("G23").Value = -3
$ws.Range("H23").Value = 13

# Row 25: fix G25 from -0 to 0
$ws.Range("G25").Value = 0

# Row 27
$ws.Range("E27").Value = 6
$ws.Range("G27").Value = -3
$ws.Range("H27").Value = 13

# Row 28: fix G28 from -0 to 0
$ws.Range("G28").Value = 0

# Match the active-cell selection left by the author (H31)
[void]$ws.Range("H31").Select()
